# Generate Report for Handoff
#
# A new handoff run (GUID 11de95f9-e9fd-40f1-bc85-4cb3bc292165) supersedes the
# previous one (68aeb95c-19ca-4db1-a4f6-7c5a94a4946d): file names/paths and
# handoff timestamps move forward, and the per-locale "Latest Target
# File"/"Latest Handback File"/"Latest Handback DateTime" columns are reset
# since this fresh handoff hasn't come back from translation yet.

$wb = $excel.ActiveWorkbook

$newGuid = "11de95f9-e9fd-40f1-bc85-4cb3bc292165"
$newFileName = "$newGuid.md"
$newPathName = "e2e\$newGuid.md"

$hyperlinkColor = 15570276   # RGB(0x64,0x95,0xED) == the workbook's custom "HyperLink" style colour

function Set-HandoffHyperlink($ws, $cellRef, $repoUrl, $displayText) {
    $range = $ws.Range($cellRef)
    $ws.Hyperlinks.Add($range, "$repoUrl/e2e/$newGuid.md", "", "", $displayText)
    $range.Font.Underline = $True
    $range.Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newFileName
$wsOverview.Range("G2").Value = "2016-08-28 13:00:39"

# Refresh the B2 hyperlink's display text (target repo/commit is unchanged).
$wsOverview.Hyperlinks.Delete()
Set-HandoffHyperlink $wsOverview "B2" "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4597fd4b7efa8ac0d1d1221ef0a8f2eac4b60411" $newPathName

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newFileName
$wsZh.Range("G2").Value = "$newGuid.5e8603efecc9af9f9779d136674a481f4ec52d1c.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-28 13:00:34"

# Latest Target File / Latest Handback File / Latest Handback DateTime aren't
# known yet for this fresh handoff.
$wsZh.Range("I2").Value = ""
$wsZh.Range("J2").Value = ""
$wsZh.Range("K2").Value = "0001-01-01 00:00:00"

# Drop the stale "Latest Target File" hyperlink that lived on I2, and
# refresh the A2 ("Source File Name") one.
$wsZh.Hyperlinks.Delete()
Set-HandoffHyperlink $wsZh "A2" "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4597fd4b7efa8ac0d1d1221ef0a8f2eac4b60411" $newFileName
$wsZh.Range("I2").Style = "Normal"

$wsZh.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsZh.Columns.Item(10).ColumnWidth = 21.7054770333426

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newFileName
$wsDe.Range("G2").Value = "$newGuid.5e8603efecc9af9f9779d136674a481f4ec52d1c.de-de.xlf"
$wsDe.Range("H2").Value = "2016-08-28 13:00:39"

# Latest Target File / Latest Handback File / Latest Handback DateTime aren't
# known yet for this fresh handoff.
$wsDe.Range("I2").Value = ""
$wsDe.Range("J2").Value = ""
$wsDe.Range("K2").Value = "0001-01-01 00:00:00"

# Drop the stale "Latest Target File" hyperlink that lived on I2, and
# refresh the A2 ("Source File Name") one.
$wsDe.Hyperlinks.Delete()
Set-HandoffHyperlink $wsDe "A2" "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4597fd4b7efa8ac0d1d1221ef0a8f2eac4b60411" $newFileName
$wsDe.Range("I2").Style = "Normal"

$wsDe.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsDe.Columns.Item(10).ColumnWidth = 21.7054770333426
